$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new price record for "Granada" at row 139,
# pushing all the existing rows 139:152 down to 140:153.
$ws.Rows.Item(139).Insert()

$ws.Cells.Item(139, 1).Value = 10
$ws.Cells.Item(139, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(139, 3).Value = 'La Araucanía'
$ws.Cells.Item(139, 4).Value = 44783
$ws.Cells.Item(139, 4).NumberFormat = $ws.Cells.Item(140, 4).NumberFormat
$ws.Cells.Item(139, 5).Value = 9
$ws.Cells.Item(139, 6).Value = 'Fruta'
$ws.Cells.Item(139, 7).Value = 100104
$ws.Cells.Item(139, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(139, 9).Value = 100104001
$ws.Cells.Item(139, 10).Value = 'Granada'
$ws.Cells.Item(139, 11).Value = 'Wonderfull'
$ws.Cells.Item(139, 12).Value = 'Primera'
$ws.Cells.Item(139, 13).Value = 125
$ws.Cells.Item(139, 14).Value = 14000
$ws.Cells.Item(139, 15).Value = 14000
$ws.Cells.Item(139, 16).Value = 14000
$ws.Cells.Item(139, 17).Value = '$/bandeja 10 kilos granel'
$ws.Cells.Item(139, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(139, 19).Value = 1400
$ws.Cells.Item(139, 20).Value = 10
